$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Row 11 ("Marking") - Right and Wrong counts
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -2

# Row 12 ("Total") - Right and Wrong totals, plus the "Right / Max" text
$ws.Range("B12").Value = 76
$ws.Range("C12").Value = -8
$ws.Range("E12").Value = "68 / 112"
